$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new test result (FAILED)
$ws.Range("A6").Value = "5a589b69-5978-4b07-9c29-72ffd3469ef7"
$ws.Range("B6").Value = "Create Citizenship by getting data from Excel"
$ws.Range("C6").Value = "FAILED"
$ws.Range("D6").Value = "2023-10-03T19:28:47.313812"
$ws.Range("E6").Value = "2023-10-03T19:28:58.174485900"
$ws.Range("F6").Value = "PT10.8606739S"

# Row 7: new test result (PASSED)
$ws.Range("A7").Value = "9e6d51ca-1be3-4fff-b576-811a84b9476c"
$ws.Range("B7").Value = "Create Citizenship by getting data from Excel"
$ws.Range("C7").Value = "PASSED"
$ws.Range("D7").Value = "2023-10-03T19:33:51.493257200"
$ws.Range("E7").Value = "2023-10-03T19:34:15.462612400"
$ws.Range("F7").Value = "PT23.9693552S"
